$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Update Only" column (M) header
$ws.Range("M1").Value = "Update Only"

# Fill "No" for every data row in the new column
$ws.Range("M2:M7").Value = "No"

# Match the formatting used by the rest of the header row / data column (style index 4)
$ws.Range("M1:M7").Style = $ws.Range("I6").Style

# Update view state to match the post-edit selection
$ws.Application.ActiveWindow.ScrollColumn = 12
$ws.Range("M3:M7").Select()
